$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-12 Sunday" "2025-10-13 Monday"

Replace-Text "950×7=6650" "845×7=5915"
Replace-Text "647×3=1941" "145×2=290"
Replace-Text "862×5=4310" "984×4=3936"
Replace-Text "900×6=5400" "839×3=2517"
Replace-Text "417×5=2085" "551×9=4959"

Replace-Text "975×4=3900" "790×3=2370"
Replace-Text "469×6=2814" "523×5=2615"
Replace-Text "225×8=1800" "645×7=4515"
Replace-Text "889×8=7112" "337×9=3033"
Replace-Text "193×6=1158" "888×8=7104"

Replace-Text "407×7=2849" "719×3=2157"
Replace-Text "666×5=3330" "593×6=3558"
Replace-Text "458×8=3664" "828×2=1656"
Replace-Text "449×7=3143" "810×9=7290"
Replace-Text "641×4=2564" "891×8=7128"

Replace-Text "775×8=6200" "930×6=5580"
Replace-Text "603×2=1206" "822×2=1644"
Replace-Text "285×9=2565" "656×4=2624"
Replace-Text "177×9=1593" "885×9=7965"
Replace-Text "848×2=1696" "159×9=1431"

Replace-Text "850×9=7650" "749×5=3745"
Replace-Text "212×2=424" "198×5=990"
Replace-Text "212×5=1060" "506×9=4554"
Replace-Text "326×9=2934" "982×2=1964"
Replace-Text "441×6=2646" "897×9=8073"
